$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the tour order (B80:B155) produced by the new Solver run (Evolutionary engine)
$ws.Range("B80").Value = 30
$ws.Range("B81").Value = 29
$ws.Range("B82").Value = 28
$ws.Range("B83").Value = 33
$ws.Range("B84").Value = 32
$ws.Range("B85").Value = 35
$ws.Range("B86").Value = 18
$ws.Range("B87").Value = 17
$ws.Range("B88").Value = 16
$ws.Range("B89").Value = 15
$ws.Range("B90").Value = 74
$ws.Range("B91").Value = 14
$ws.Range("B92").Value = 13
$ws.Range("B93").Value = 12
$ws.Range("B94").Value = 9
$ws.Range("B95").Value = 6
$ws.Range("B96").Value = 5
$ws.Range("B97").Value = 4
$ws.Range("B98").Value = 3
$ws.Range("B99").Value = 7
$ws.Range("B100").Value = 8
$ws.Range("B101").Value = 2
$ws.Range("B102").Value = 76
$ws.Range("B103").Value = 75
$ws.Range("B104").Value = 1
$ws.Range("B105").Value = 23
$ws.Range("B106").Value = 22
$ws.Range("B107").Value = 24
$ws.Range("B108").Value = 21
$ws.Range("B109").Value = 25
$ws.Range("B110").Value = 26
$ws.Range("B111").Value = 27
$ws.Range("B112").Value = 46
$ws.Range("B113").Value = 45
$ws.Range("B114").Value = 44
$ws.Range("B115").Value = 48
$ws.Range("B116").Value = 47
$ws.Range("B117").Value = 69
$ws.Range("B118").Value = 68
$ws.Range("B119").Value = 70
$ws.Range("B120").Value = 67
$ws.Range("B121").Value = 50
$ws.Range("B122").Value = 49
$ws.Range("B123").Value = 51.000000000000007
$ws.Range("B124").Value = 66
$ws.Range("B125").Value = 65
$ws.Range("B126").Value = 71
$ws.Range("B127").Value = 72
$ws.Range("B128").Value = 73
$ws.Range("B129").Value = 63.999999999999993
$ws.Range("B130").Value = 63
$ws.Range("B131").Value = 61.999999999999993
$ws.Range("B132").Value = 61
$ws.Range("B133").Value = 60
$ws.Range("B134").Value = 59
$ws.Range("B135").Value = 57.999999999999993
$ws.Range("B136").Value = 57
$ws.Range("B137").Value = 56
$ws.Range("B138").Value = 55
$ws.Range("B139").Value = 52
$ws.Range("B140").Value = 53
$ws.Range("B141").Value = 54
$ws.Range("B142").Value = 42
$ws.Range("B143").Value = 43
$ws.Range("B144").Value = 34
$ws.Range("B145").Value = 41
$ws.Range("B146").Value = 40
$ws.Range("B147").Value = 39
$ws.Range("B148").Value = 38
$ws.Range("B149").Value = 36
$ws.Range("B150").Value = 37
$ws.Range("B151").Value = 11
$ws.Range("B152").Value = 10
$ws.Range("B153").Value = 20
$ws.Range("B154").Value = 19
$ws.Range("B155").Value = 31

# Switch Solver engine setting from GRG Nonlinear (2) to Evolutionary (3)
$wb.Names("solver_eng").RefersTo = "=3"

# Update the active selection to reflect the last-used cell after recalculation
[void]$ws.Range("E156").Select()
